$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row for the added columns F..K
$ws.Cells.Item(1, 6).Value  = "frequency"
$ws.Cells.Item(1, 7).Value  = "frequency_occurrence"
$ws.Cells.Item(1, 8).Value  = "frequency_occurrence_probab"
$ws.Cells.Item(1, 9).Value  = "max_probab"
$ws.Cells.Item(1, 10).Value = "max_probab_percentage"
$ws.Cells.Item(1, 11).Value = "recommended_level"

# Per-row data: frequency, frequency_occurrence (json), frequency_occurrence_probab (json),
# max_probab, max_probab_percentage (text), recommended_level (moved from old column F)
$data = @{
    2  = @(8,  '{"L3":6,"L2":2}',             '{"L3":0.75,"L2":0.25}',               0.75,  "75.00",  'L3')
    3  = @(8,  '{"L3":4,"L2":3,"L1":1}',      '{"L3":0.5,"L2":0.375,"L1":0.125}',    0.5,   "50.00",  'L3')
    4  = @(8,  '{"L3":7,"L2":1}',             '{"L3":0.875,"L2":0.125}',             0.875, "87.50",  'L3')
    5  = @(8,  '{"L2":7,"L1":1}',             '{"L2":0.875,"L1":0.125}',             0.875, "87.50",  'L2')
    6  = @(11, '{"L3":11}',                   '{"L3":1.0}',                          1,     "100.00", 'L3')
    7  = @(11, '{"L3":11}',                   '{"L3":1.0}',                          1,     "100.00", 'L3')
    8  = @(8,  '{"L3":7,"L2":1}',             '{"L3":0.875,"L2":0.125}',             0.875, "87.50",  'L3')
    9  = @(8,  '{"L3":7,"L2":1}',             '{"L3":0.875,"L2":0.125}',             0.875, "87.50",  'L3')
    11 = @(8,  '{"L2":6,"L1":2}',             '{"L2":0.75,"L1":0.25}',               0.75,  "75.00",  'L2')
    14 = @(8,  '{"L2":8}',                    '{"L2":1.0}',                          1,     "100.00", 'L2')
    15 = @(8,  '{"L2":7,"L3":1}',             '{"L2":0.875,"L3":0.125}',             0.875, "87.50",  'L3')
    16 = @(8,  '{"L3":6,"L2":2}',             '{"L3":0.75,"L2":0.25}',               0.75,  "75.00",  'L3')
    17 = @(8,  '{"L3":8}',                    '{"L3":1.0}',                          1,     "100.00", 'L3')
    18 = @(8,  '{"L3":8}',                    '{"L3":1.0}',                          1,     "100.00", 'L3')
    19 = @(8,  '{"L3":8}',                    '{"L3":1.0}',                          1,     "100.00", 'L3')
    20 = @(8,  '{"L3":7,"L1":1}',             '{"L3":0.875,"L1":0.125}',             0.875, "87.50",  'L3')
    21 = @(8,  '{"L2":5,"L3":3}',             '{"L2":0.625,"L3":0.375}',             0.625, "62.50",  'L3')
    22 = @(8,  '{"L2":5,"L1":3}',             '{"L2":0.625,"L1":0.375}',             0.625, "62.50",  'L2')
    23 = @(8,  '{"L3":7,"L2":1}',             '{"L3":0.875,"L2":0.125}',             0.875, "87.50",  'L3')
    24 = @(8,  '{"L3":8}',                    '{"L3":1.0}',                          1,     "100.00", 'L3')
    25 = @(8,  '{"L3":8}',                    '{"L3":1.0}',                          1,     "100.00", 'L3')
    26 = @(3,  '{"L2":3}',                    '{"L2":1.0}',                          1,     "100.00", 'L2')
    30 = @(3,  '{"L2":3}',                    '{"L2":1.0}',                          1,     "100.00", 'L2')
    41 = @(8,  '{"L3":5,"L2":2,"L1":1}',      '{"L3":0.625,"L2":0.25,"L1":0.125}',   0.625, "62.50",  'L3')
}

# Rows that stay blank across F..K (same rows whose old column F was already blank)
$emptyRows = @(10, 12, 13, 27, 28, 29, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 42, 43)

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 6).Value  = $vals[0]
    $ws.Cells.Item($row, 7).Value  = $vals[1]
    $ws.Cells.Item($row, 8).Value  = $vals[2]
    $ws.Cells.Item($row, 9).Value  = $vals[3]
    # Force text storage for the percentage string (otherwise "75.00" auto-coerces to the number 75)
    $ws.Cells.Item($row, 10).Value = "'" + $vals[4]
    $ws.Cells.Item($row, 11).Value = $vals[5]
}

foreach ($row in $emptyRows) {
    for ($col = 6; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Style = "Normal"
    }
}
